# Populate the "header4" report data into Sheet1.
# Before: an empty templated sheet with only a styled header row (row 8, A8:F8).
# After:  the generated report — title/meta rows, the styled header row with
#         its labels, and one data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: report title + generation timestamp
$ws.Cells.Item(1, 1).Value = "KAMOTO COPPER COMPANY OFFLINE DEVICES"
$ws.Cells.Item(1, 3).Value = "2024-01-11 07:50:32"

# Row 2 stays blank

# Rows 3-6: report meta lines
$ws.Cells.Item(3, 1).Value = "Clients: KAMOTO COPPER COMPANY"
$ws.Cells.Item(4, 1).Value = "Max Hours: 96"
$ws.Cells.Item(5, 1).Value = "Schedule: VEHICLES OFFLINE REPORT"
$ws.Cells.Item(6, 1).Value = "KAMOTO COPPER COMPANY"

# Row 7 stays blank

# Row 8: column headers (already carries the "Good" cell style, s="1",
# from the template - just fill in the text)
$ws.Cells.Item(8, 1).Value = "Device"
$ws.Cells.Item(8, 2).Value = "Asset"
$ws.Cells.Item(8, 3).Value = "Last Date"
$ws.Cells.Item(8, 4).Value = "Last Update"
$ws.Cells.Item(8, 5).Value = "Location"
$ws.Cells.Item(8, 6).Value = "Remarque"

# Row 9: the single offline-device data row
$ws.Cells.Item(9, 1).Value = 1011
# Asset is a text value that happens to look numeric - force text (via the
# leading apostrophe) so it isn't coerced to a number, then drop back to the
# plain "Normal" style so the quote-prefix marker doesn't linger as a style.
$assetCell = $ws.Cells.Item(9, 2)
$assetCell.Value = "'131"
$assetCell.Style = "Normal"
$ws.Cells.Item(9, 3).Value = "2024-01-10 12:57:58"
$ws.Cells.Item(9, 4).Value = "1 day ago"
$ws.Cells.Item(9, 5).Value = "KCC, Luilu, Lualaba, Congo - Kinshasa"
$ws.Cells.Item(9, 6).Value = "0:02:48"
